$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, D(open_price), E(close_price), F(high_price), G(low_price), H(shares_outstanding)
$data = @(
    @(2, 89.26210374595165, 85.79873657226562, 92.04708383448224, 83.41543322652859, 252114595),
    @(3, 89.6816417642498, 86.5753173828125, 90.30647380153049, 76.51548072164294, 252114595),
    @(4, 76.84572750590219, 69.93684387207031, 87.5571779565623, 64.74179067590643, 252114595),
    @(5, 73.19492317756635, 66.75019836425781, 76.55118009840629, 59.74312476355895, 252114595),
    @(6, 71.59712337694421, 76.12271118164062, 81.38917598645428, 71.45430757062508, 252114595),
    @(7, 69.47268408575884, 75.06048583984375, 79.53251831434947, 66.98227103624932, 252114595),
    @(8, 91.92212165949459, 89.26210784912109, 95.99247159816385, 87.36082740359183, 252114595),
    @(9, 87.70002159777536, 87.34297180175781, 88.18203780087526, 85.69162415660493, 252114595),
    @(10, 92.54695090405528, 94.39467620849609, 94.4393099851402, 92.08279095368404, 252114595),
    @(11, 97.74201769510552, 98.48289489746094, 98.57215565204568, 97.29570711202284, 252114595),
    @(12, 101.0000841479672, 104.4813079833984, 105.3292953490091, 100.732295079796, 252114595),
    @(13, 104.4366750423875, 107.4001770019531, 108.417764572938, 104.2581535356094, 252114595),
    @(14, 104.4366443523244, 93.63593292236328, 104.757988418221, 89.46739499162697, 252114595),
    @(15, 96.938651101509, 85.10249328613281, 99.51832556240194, 80.55905510897375, 252114595),
    @(16, 77.3555743563485, 67.11667633056641, 77.85677686277239, 62.3999851846563, 252114595),
    @(17, 64.56246823987456, 78.13755035400391, 79.00843993563362, 64.24822464661639, 252114595),
    @(18, 81.02225649341855, 95.08412170410156, 96.95663085645415, 81.02225649341855, 252114595),
    @(19, 91.98946213367324, 93.32538604736328, 95.91600521617022, 84.80433103882154, 252114595),
    @(20, 101.1168953064144, 102.9646072387695, 107.3302773138598, 92.28592312236268, 252114595),
    @(21, 116.8407326101147, 115.2508773803711, 125.2987602143609, 113.2703717312131, 252114595),
    @(22, 72.50229647033137, 90.82855224609376, 98.36339243774692, 66.46347687426717, 252114595),
    @(23, 104.4194807235422, 107.6056213378906, 115.4702748203638, 102.2496039196007, 252114595),
    @(24, 119.041726486157, 124.0832901000977, 133.294030872814, 116.2959469389671, 252114595),
    @(25, 147.5863992591716, 147.7152709960938, 168.0402752290201, 143.6189757375344, 252114595),
    @(26, 189.5204891505387, 177.739501953125, 199.8242179812948, 176.3915295803338, 252114595),
    @(27, 190.3024250274495, 191.0894165039062, 193.2744600666496, 171.2110792892394, 252114595),
    @(28, 182.2529814664118, 186.4581756591797, 188.8438972091272, 169.266070713695, 252114595),
    @(29, 212.5538281334956, 191.1951446533203, 218.6124307748278, 167.6307832560902, 252114595),
    @(30, 174.9734117405344, 159.8063049316406, 174.9734117405344, 154.0555276639445, 252114595),
    @(31, 137.8403395505982, 172.8102722167969, 174.1165924611291, 131.8820163837993, 252114595),
    @(32, 141.2867433746142, 138.0361633300781, 154.7709609918054, 124.8070678865, 252114595),
    @(33, 151.7172528416604, 175.0094299316406, 175.807038754744, 145.3743429310661, 252114595),
    @(34, 177.1949963306114, 156.3839263916016, 177.6438825761018, 148.3326396193198, 252114595),
    @(35, 197.3629480785773, 214.109130859375, 216.5961034242072, 192.0721589607045, 252114595),
    @(36, 193.7577806985518, 166.3993072509766, 198.9303169527541, 161.3618886763812, 252114595),
    @(37, 220.324586579236, 204.1439361572266, 220.324586579236, 195.4282831980479, 252114595),
    @(38, 241.8374958874331, 249.4417724609375, 254.2711141532573, 207.4089032709698, 252114595),
    @(39, 262.9734662228032, 257.1495971679688, 289.3177086830925, 239.8147620758475, 252114595),
    @(40, 234.3851839165721, 230.1550445556641, 251.865187249768, 223.1178963763991, 252114595),
    @(41, 206.6631526022086, 205.6180572509766, 216.8183508097865, 199.4066241691167, 252114595),
    @(42, 186.8138951822417, 182.6129302978516, 196.3254923945998, 146.7264312152305, 252114595),
    @(43, 215.5961915656884, 212.7891387939453, 235.3750300894592, 211.2163866668767, 252114595),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = $row[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "NXPI"    # I: fixed_ticker
}
